$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.929.16"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.811.34"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'595.44"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").Value = "'174.38"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").Value = "3.809.99"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").Value = "'38.02"
$ws.Range("E13").Value = "  -4.63%  "
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "4.447.05"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "3.817.72"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "68.059.95"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").Value = "'7.14"
$ws.Range("E19").Value = "  -5.28%  "
$ws.Range("D20").Value = "'16.22"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "'489.46"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "'9.19"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "'84.60"
$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  -8.87%  "
$ws.Range("D27").Value = "'12.30"
$ws.Range("E27").Value = "  -5.41%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  -9.52%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'2.42"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "'32.71"
$ws.Range("E32").Value = "  +6.82%  "
$ws.Range("D33").Value = "'7.71"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("D38").Value = "'5.77"
$ws.Range("E38").Value = "  -6.26%  "
$ws.Range("D39").Value = "'0.325"
$ws.Range("E39").Value = "  -7.06%  "
$ws.Range("D40").Value = "'448.31"
$ws.Range("E40").Value = "  +2.45%  "
$ws.Range("D41").Value = "'48.89"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").Value = "'2.89"
$ws.Range("E43").Value = "  -6.09%  "
$ws.Range("D44").Value = "'8.26"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "'41.47"
$ws.Range("E45").Value = "  -9.10%  "
$ws.Range("D46").Value = "2.832.57"
$ws.Range("E46").Value = "  -4.91%  "
$ws.Range("D48").Value = "'138.61"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  -3.87%  "
$ws.Range("D50").Value = "'26.12"
$ws.Range("E50").Value = "  -4.65%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "  -6.60%  "
